$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DIGIKEY part numbers (column H) that changed due to revised
# part/order codes.
$ws.Range("H3").Value  = "RHM1.0KCECT-ND"

# Updated Part Name for the LTC1966 RMS-to-DC converter row (U5, was IC1)
$ws.Range("F15").Value = "U5"

$ws.Range("H5").Value  = "RHM14KCDCT-ND"
$ws.Range("H6").Value  = "RHM20.0KCDCT-ND"
$ws.Range("H7").Value  = "RHM47KCECT-ND"
$ws.Range("H8").Value  = "490-6318-1-ND"
$ws.Range("H9").Value  = "1276-1442-1-ND"
$ws.Range("H10").Value = "490-6281-1-ND"
$ws.Range("H13").Value = "497-6871-1-ND"
$ws.Range("H16").Value = "FK3503010LCT-ND"

# Highlight the header/part row (row 2) in yellow to call it out
$ws.Range("A2:H2").Interior.Color = 65535

# Leave the selection where the editor last clicked
$ws.Range("H17").Select() | Out-Null
